# Generate Report for Handoff
# Updates the "Priority" column (E) for rows that are ready for handoff
# (rows 7, 8, 9, 10, 12, 13 on both the zh-cn and de-de sheets) from blank
# to "ht", and refreshes the "Latest Handoff Datetime" / "Latest Handback
# DateTime" timestamps for those same rows (column H on the language
# sheets, which rolls up into column G "Latest HO Xliff Generate Date" on
# the Overview sheet).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 13)

# --- zh-cn sheet: Priority -> "ht", handoff timestamp advanced by 14s ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("E$r").Value = "ht"
    $wsZh.Range("H$r").Value = "2016-08-27 00:22:43"
}

# --- de-de sheet: Priority -> "ht", handoff timestamp advanced by 14s ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("E$r").Value = "ht"
    $wsDe.Range("H$r").Value = "2016-08-27 00:22:48"
}

# --- Overview sheet: rollup of the latest handoff/handback datetime ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-27 00:22:48"
}
